$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "37.669.93"
Set-TextValue $ws "E2" "  +0.85%  "
Set-TextValue $ws "D3" "2.073.05"
Set-TextValue $ws "E3" "  +0.29%  "
Set-TextValue $ws "E4" "  +0.03%  "
Set-TextValue $ws "D5" "232.34"
Set-TextValue $ws "E5" "  -0.56%  "
Set-TextValue $ws "D6" "0.624"
Set-TextValue $ws "E6" "  +0.56%  "
Set-TextValue $ws "E7" "  -0.04%  "
Set-TextValue $ws "D8" "56.84"
Set-TextValue $ws "E8" "  +0.14%  "
Set-TextValue $ws "E9" "  +0.82%  "
Set-TextValue $ws "D10" "0.0785"
Set-TextValue $ws "E10" "  +2.98%  "
Set-TextValue $ws "E11" "  +2.76%  "
Set-TextValue $ws "D12" "2.366.15"
Set-TextValue $ws "E12" "  -0.28%  "
Set-TextValue $ws "D13" "14.37"
Set-TextValue $ws "E13" "  -0.27%  "
Set-TextValue $ws "D14" "20.81"
Set-TextValue $ws "E14" "  +0.28%  "
Set-TextValue $ws "E15" "  -2.46%  "
Set-TextValue $ws "E16" "  +1.60%  "
Set-TextValue $ws "D17" "2.073.93"
Set-TextValue $ws "E17" "  +0.26%  "
Set-TextValue $ws "D18" "37.592.66"
Set-TextValue $ws "E18" "  +0.77%  "
Set-TextValue $ws "E19" "  -4.12%  "
Set-TextValue $ws "D20" "70.66"
Set-TextValue $ws "E20" "  +1.66%  "
Set-TextValue $ws "E21" "  +0.47%  "
Set-TextValue $ws "D22" "227.29"
Set-TextValue $ws "E22" "  +0.60%  "
Set-TextValue $ws "D23" "1.00"
Set-TextValue $ws "E23" "  +0.06%  "
Set-TextValue $ws "E24" "  -1.14%  "
Set-TextValue $ws "E25" "  -0.93%  "
Set-TextValue $ws "D26" "169.71"
Set-TextValue $ws "E26" "  +2.01%  "
Set-TextValue $ws "E27" "  +10.31%  "
Set-TextValue $ws "D28" "8.86"
Set-TextValue $ws "E28" "  +0.94%  "
Set-TextValue $ws "E29" "  -0.92%  "
Set-TextValue $ws "E30" "  +1.92%  "
Set-TextValue $ws "E31" "  +0.92%  "
Set-TextValue $ws "D32" "4.60"
Set-TextValue $ws "E32" "  +2.09%  "
Set-TextValue $ws "E33" "  +0.46%  "
Set-TextValue $ws "E34" "  -0.55%  "
Set-TextValue $ws "E35" "  -0.07%  "
Set-TextValue $ws "E36" "  +3.62%  "
Set-TextValue $ws "D37" "3.35"
Set-TextValue $ws "E37" "  +4.05%  "
Set-TextValue $ws "E38" "  +0.14%  "
Set-TextValue $ws "E39" "  -4.68%  "
Set-TextValue $ws "D40" "0.0990"
Set-TextValue $ws "E40" "  +6.12%  "
Set-TextValue $ws "B41" "Aave"
Set-TextValue $ws "C41" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D41" "98.94"
Set-TextValue $ws "E41" "  +3.03%  "
Set-TextValue $ws "B42" "HuobiToken"
Set-TextValue $ws "C42" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws "D42" "2.93"
Set-TextValue $ws "E42" "  -0.79%  "
Set-TextValue $ws "B43" "FTXToken"
Set-TextValue $ws "C43" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws "D43" "4.39"
Set-TextValue $ws "E43" "  +3.85%  "
Set-TextValue $ws "B44" "VeChain"
Set-TextValue $ws "C44" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D44" "0.0213"
Set-TextValue $ws "E44" "  +0.80%  "
Set-TextValue $ws "D45" "1.451.13"
Set-TextValue $ws "E45" "  -1.60%  "
Set-TextValue $ws "E46" "  -1.65%  "
Set-TextValue $ws "E47" "  +2.36%  "
Set-TextValue $ws "B48" "InjectiveProtocol"
Set-TextValue $ws "C48" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws "D48" "15.50"
Set-TextValue $ws "E48" "  +1.65%  "
Set-TextValue $ws "B49" "FraxShare"
Set-TextValue $ws "C49" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws "D49" "7.38"
Set-TextValue $ws "E49" "  +2.70%  "
Set-TextValue $ws "E50" "  +1.09%  "
Set-TextValue $ws "D51" "47.41"
Set-TextValue $ws "E51" "  +8.00%  "
